$d = $word.ActiveDocument

# The title paragraph currently reads "Clase 1 – ReactJS".
# Replace the leading "Clase " with "Módulo " (new module 3 naming).
$rng = $d.Paragraphs(1).Range
$rng.Find.Execute("Clase ", $true, $false, $false, $false, $false, $true, 1, $false, "Módulo ", 2)
